# Rename 'Codelists' sheet to 'Cells' (Close #256)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Codelists")
$ws.Name = "Cells"

# The renamed sheet becomes the active/selected tab (was previously on
# "Table", the first sheet).
$ws.Activate()

# Update the active sheet's selection to reflect the new focus.
$ws.Range("F19").Select()
